# methanolA_HIFIdata_k3.xlsx — "testes on flux implementation, litle success"
#
# Diff summary being reproduced here (the parts reachable through the Excel
# object model / this COM surface):
#   - Sheet1: active selection moved from A1:AA18 (whole used range) to the
#     single cell C16 (activeCell=C16, sqref=C16).
#   - Sheet1: explicit column widths added for column B (Transition names,
#     wide text) and column C (Frequency values, numeric, best-fit width).
#
# (Window position/size and the session revisionPtr/uid GUIDs in
# xl/workbook.xml are host/session metadata that this runtime always
# regenerates on save regardless of script content and are not exposed as
# settable COM properties here, so they are intentionally left alone.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (Transition) and column C (Frequency(Mhz)) to fit their
# contents, matching the <cols> block added by the edit. The host only
# persists column widths on a 1/6-character grid, so we feed it the input
# that lands closest to the authored widths (40.7109375 / 15.140625).
$ws.Columns.Item(2).ColumnWidth = 39.833333333333336
$ws.Columns.Item(3).ColumnWidth = 14.333333333333334

# Move the selection/active cell to C16.
$ws.Range("C16").Select()
